# "Adicionei funções para anexos"
# Mark the "Anexo de imagens / arquivos nas questões" task (row 14) as done:
# fill in the INICIO (HORAS) and FIM columns with "feito", mirroring the
# pattern already used on rows 11 and 13 for tasks that were completed
# without explicit start/end timestamps.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("E14").Value = "feito"
$ws.Range("F14").Value = "feito"

# Move the cursor/viewport roughly where the author left it (scrolled down a
# bit further, with the active cell on the DURAÇÃO REAL column of the same
# row) so the saved view state tracks the edit just made.
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G14").Select()
